# Update the localization Status for two files from "Ready for handoff"
# to "In Translation" across the Overview sheet (zh-cn/de-de columns) and
# the per-locale "zh-cn" / "de-de" report sheets.
#
# Affected files (rows 3 and 4 in every sheet):
#   577fe84a-60dd-450f-90fd-a4f4812983db.md
#   df39a1bb-6da6-4d53-9e82-f60e118a574e.md
#
# Row 5 (f54534a4-e901-482a-bef5-d6fd27abb9cf.md) keeps its
# "Ready for handoff" status and is left untouched.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

# --- zh-cn sheet: column C holds the status ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

# --- de-de sheet: column C holds the status ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus
